$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.908.51"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.649.86"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'587.95"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").Value = "'144.51"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'6.58"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").Value = "'0.381"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "3.119.00"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "'26.06"
$ws.Range("E14").Value = "  +11.19%  "
$ws.Range("D15").Value = "60.950.66"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "'0.0000144"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "2.662.67"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("E18").Value = "  +2.66%  "
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "'350.24"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'6.91"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'0.531"
$ws.Range("D24").Value = "'64.04"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'8.14"
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("D28").Value = "'1.99"
$ws.Range("E28").Value = "  +8.90%  "
$ws.Range("D29").Value = "0.0₃0813"
$ws.Range("E29").Value = "  +2.06%  "
$ws.Range("D30").Value = "'6.86"
$ws.Range("E30").Value = "  +7.78%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'164.64"
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("D33").Value = "'19.93"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  +7.08%  "
$ws.Range("E35").Value = "  +1.87%  "
$ws.Range("E36").Value = "  +7.57%  "
$ws.Range("D37").Value = "'339.95"
$ws.Range("E37").Value = "  +13.27%  "
$ws.Range("E38").Value = "  +3.14%  "
$ws.Range("D39").Value = "'4.08"
$ws.Range("E39").Value = "  +4.70%  "
$ws.Range("D40").Value = "'0.911"
$ws.Range("E40").Value = "  +7.26%  "
$ws.Range("D41").Value = "'38.62"
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("D42").Value = "'5.27"
$ws.Range("E42").Value = "  +5.50%  "
$ws.Range("D43").Value = "'20.42"
$ws.Range("E43").Value = "  +2.25%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.622"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0250"
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'133.27"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'20.65"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0562"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("D49").Value = "'0.0996"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "2.108.57"
$ws.Range("E51").Value = "  +4.03%  "
